$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mops56's UTC value changed from 16 to 18
$ws.Range("D22").Value = 18

# Re-sort the Table1 data by the UTC column (column D), ascending,
# matching the table's existing sortState condition.
$lo = $ws.ListObjects.Item(1)
$lo.Sort.SortFields.Clear()
$lo.Sort.SortFields.Add($ws.Range("D2:D40"))
$lo.Sort.Header = 1
$lo.Sort.Apply()

# Selected cell moved to D1 (the UTC header, used to trigger the sort)
$null = $ws.Range("D1").Select()
